$d = $word.ActiveDocument

function Insert-ParagraphXml($para, $innerBody) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        $innerBody +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml)
}

# Insert 4 blank placeholder paragraphs before the existing first paragraph
# ("Phase C - Logical Design"), which will become paragraph 5.
$firstRange = $d.Paragraphs(1).Range
$firstRange.InsertParagraphBefore()
$firstRange.InsertParagraphBefore()
$firstRange.InsertParagraphBefore()
$firstRange.InsertParagraphBefore()

# Paragraph 1: title, centered
$titlePara = $d.Paragraphs(1)
$titleBody = '<w:body><w:p>' +
    '<w:pPr><w:pStyle w:val="segoe"/><w:ind w:left="0"/><w:jc w:val="center"/>' +
    '<w:rPr><w:rFonts w:cs="Segoe UI" w:hint="eastAsia"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cs="Segoe UI" w:hint="eastAsia"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr>' +
    '<w:t>CMSC 461 Databases Final Project Phases C, D, E, F</w:t></w:r>' +
    '</w:p></w:body>'
Insert-ParagraphXml $titlePara $titleBody

# Paragraph 2: author name, centered
$namePara = $d.Paragraphs(2)
$nameBody = '<w:body><w:p>' +
    '<w:pPr><w:pStyle w:val="segoe"/><w:ind w:left="0"/><w:jc w:val="center"/>' +
    '<w:rPr><w:rFonts w:cs="Segoe UI" w:hint="eastAsia"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cs="Segoe UI" w:hint="eastAsia"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr>' +
    '<w:t>Khang Ngo</w:t></w:r>' +
    '</w:p></w:body>'
Insert-ParagraphXml $namePara $nameBody

# Paragraph 3: author email, centered
$emailPara = $d.Paragraphs(3)
$emailBody = '<w:body><w:p>' +
    '<w:pPr><w:pStyle w:val="segoe"/><w:ind w:left="0"/><w:jc w:val="center"/>' +
    '<w:rPr><w:rFonts w:cs="Segoe UI" w:hint="eastAsia"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cs="Segoe UI" w:hint="eastAsia"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr>' +
    '<w:t>khang4@umbc.edu</w:t></w:r>' +
    '</w:p></w:body>'
Insert-ParagraphXml $emailPara $emailBody

# Paragraph 4: blank bold paragraph (no run)
$blankPara = $d.Paragraphs(4)
$blankBody = '<w:body><w:p>' +
    '<w:pPr><w:pStyle w:val="segoe"/><w:ind w:left="0"/>' +
    '<w:rPr><w:rFonts w:cs="Segoe UI" w:hint="eastAsia"/><w:b/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr>' +
    '</w:pPr>' +
    '</w:p></w:body>'
Insert-ParagraphXml $blankPara $blankBody

Write-Output ("paragraph count: " + $d.Paragraphs.Count)
